$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume/change (column E) values.
# Values are stored as text, so number format is forced to "@" (Text) before
# assignment to avoid Excel auto-converting them to numeric/percentage values,
# which would alter their textual representation (e.g. trailing zeros, "%" sign).
$cellUpdates = @{
    "D2" = "265.80"
    "E2" = "1.84%"
    "D3" = "26.59"
    "E3" = "-1.77%"
    "D4" = "4.699"
    "E4" = "-0.10%"
    "D5" = "0.06085"
    "D6" = "6.737"
    "E6" = "0.77%"
    "D7" = "0.8504"
    "E7" = "-0.07%"
    "D8" = "0.9108"
    "E8" = "-0.54%"
    "D9" = "0.1409"
    "E9" = "0.17%"
    "E10" = "10.24%"
    "D11" = "0.07085"
    "E11" = "-0.07%"
    "D12" = "0.03150"
    "E12" = "-0.19%"
    "D13" = "0.09030"
    "D14" = "0.001527"
    "E14" = "-0.17%"
    "D15" = "0.0006070"
    "E15" = "-1.50%"
    "D16" = "0.005999"
    "E16" = "-0.95%"
    "D17" = "3.455"
    "E17" = "0.09%"
    "E18" = "-0.38%"
    "E19" = "0.31%"
    "D22" = "4.085"
    "E22" = "0.16%"
    "D23" = "0.04260"
    "E23" = "1.17%"
    "E24" = "-3.11%"
    "E25" = "6.83%"
    "E26" = "0.03%"
    "D40" = "0.03922"
    "E40" = "0.90%"
    "D41" = "0.1114"
    "E41" = "0.32%"
    "D42" = "0.004188"
    "E42" = "2.14%"
    "D44" = "0.01161"
    "E44" = "-28.87%"
    "D45" = "0.00005108"
    "E45" = "-1.05%"
    "E48" = "54.93%"
}

foreach ($cellRef in $cellUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellUpdates[$cellRef]
}

